# Dataframe ST.xlsx update: replace the "25-sep".."29-sep" (CB:CE) block with a
# single new "25 - Oct" column of fresher data, and refresh the VLOOKUP helper
# table on Sheet3 so every product resolves (no more #N/A).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# 1) Sheet3: fix the three DORITOS labels (typo "GR" -> "G") in the lookup
#    table (A20:A36), and load each product's new "25 - Oct" figure (B20:B36).
# ---------------------------------------------------------------------------

$lookupNames = @(
    "3D QUESO 92GX27",
    "CHEETOS 94GRX24",
    "DORITOS QUESO 140GX19",
    "DORITOS QUESO 40GX58X1 CH",
    "DORITOS QUESO 85GX26",
    "LAYS CEBOLLA CARAMELIZADA 85GX25",
    "LAYS CLASICAS 145GRX18",
    "LAYS CLASICAS 249GRX14",
    "LAYS CLASICAS 40GX68",
    "LAYS CLASICAS 94GRX25",
    "LAYS ONDAS FH 30GX72",
    "LAYS ONDAS FH 70GX28",
    "LAYS QSO Y CEBOLLA 34GX72",
    "PEHUAMAR ACANALADA 520GX9",
    "PEHUAMAR MAICITOS 285GX10",
    "PEHUAMAR PAPA LISA 520GX9",
    "QUAKER AVENA INSTANT FORTIF 18X280G"
)

$lookupValues = @(
    19.157121981040763,
    2.1926697864742222,
    18.233398417051458,
    0,
    0,
    0,
    5.8273851245061286,
    9.0340270253750035,
    4.6471908629674861,
    6.1492215874504037,
    4.5005648034463546,
    5.3926617198736801,
    5.1325978829223429,
    6.9594126652105421,
    5.5233505224672443,
    7.6196685419009063,
    29.777407690159443
)

for ($i = 0; $i -lt $lookupNames.Length; $i++) {
    $row = 20 + $i
    $ws3.Cells.Item($row, 1).Value = $lookupNames[$i]
    $ws3.Cells.Item($row, 2).Value = $lookupValues[$i]
}

# ---------------------------------------------------------------------------
# 2) Sheet3: B2:B18 now references the same (corrected) product names, in the
#    same order as the lookup table, so every VLOOKUP in C2:C18 resolves.
# ---------------------------------------------------------------------------

for ($i = 0; $i -lt $lookupNames.Length; $i++) {
    $row = 2 + $i
    $ws3.Cells.Item($row, 2).Value = $lookupNames[$i]
}

# Column A width tightened slightly on Sheet3.
$ws3.Columns.Item(1).ColumnWidth = 36.28515625

# ---------------------------------------------------------------------------
# 3) Sheet1: the CA column becomes the new "25 - Oct" figures (same row order
#    as the Sheet3 lookup table), then CB:CE (the old 26/27/28/29-sep columns)
#    are deleted outright, shifting everything left.
# ---------------------------------------------------------------------------

$ws1.Range("CA1").Value = "25 - Oct"

for ($i = 0; $i -lt $lookupValues.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 79).Value = $lookupValues[$i]   # column 79 = CA
}

$ws1.Range("CB1:CE18").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 4) Misc view state that travelled with the edit.
# ---------------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("BX6").Select()
